# Fix formatting on fastq purpose column: "fullRNASEQ" -> "fullRNASeq"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
